$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Environment")
$ws2 = $wb.Worksheets.Item("Aragonite")

# Add the two new rows of environmental parameter data to the Environment
# sheet (current velocity + wave height limits), in the same order the
# values were originally typed so the shared-string table lines up.
$ws1.Range("B8").Value = "Current velocity"
$ws1.Range("B9").Value = "Wave height"
$ws1.Range("C8").Value = "0.04 - 1.6 m/s"
$ws1.Range("D8").Value = "0.04 - 1.6 m/s"
$ws1.Range("E8").Value = "Froehlich et al. (2017)"
$ws1.Range("C9").Value = "< 12 m"
$ws1.Range("E9").Value = "Froehlich et al. (2017)"
$ws1.Range("A8").Value = "uo / vo (m/s)"

# Move the active sheet/selection: Environment becomes the selected tab
# with F3 selected, and Aragonite is no longer the selected tab.
$ws1.Activate()
$ws1.Range("F3").Select()

# Reposition the saved window location.
$win = $wb.Windows.Item(1)
$win.Left = 17520
$win.Top = 4600
